# Update Overview Decks for March 2025
# Change the table style applied to every table in the deck from the old
# custom style GUID to the new style GUID.

$oldStyleId = "{591A0ED6-D4D7-4D79-807F-AB3D0AB6F1A9}"
$newStyleId = "{F126362D-DFA5-4DE9-A93B-187AFD16EC08}"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
